# Applies the roster update:
#  - Column A dates switch from text strings to real Excel dates,
#    formatted with a custom number format "YYYY-MM-DD HH:MM:SS".
#  - Column D (Points) bumps from 5 to 6 for every "Purple" turn.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Year = 2018; Month = 2; Day = 10; Timetable = "Purple"; Points = 6 },
    @{ Row = 3;  Year = 2018; Month = 2; Day = 11; Timetable = "Yellow"; Points = 3 },
    @{ Row = 4;  Year = 2018; Month = 2; Day = 11; Timetable = "Yellow"; Points = 4 },
    @{ Row = 5;  Year = 2018; Month = 2; Day = 11; Timetable = "Yellow"; Points = 2 },
    @{ Row = 6;  Year = 2018; Month = 2; Day = 14; Timetable = "Purple"; Points = 6 },
    @{ Row = 7;  Year = 2018; Month = 2; Day = 15; Timetable = "Purple"; Points = 6 },
    @{ Row = 8;  Year = 2018; Month = 2; Day = 17; Timetable = "Purple"; Points = 6 },
    @{ Row = 9;  Year = 2018; Month = 2; Day = 18; Timetable = "Yellow"; Points = 3 },
    @{ Row = 10; Year = 2018; Month = 2; Day = 18; Timetable = "Yellow"; Points = 4 },
    @{ Row = 11; Year = 2018; Month = 2; Day = 18; Timetable = "Yellow"; Points = 2 },
    @{ Row = 12; Year = 2018; Month = 2; Day = 21; Timetable = "Purple"; Points = 6 },
    @{ Row = 13; Year = 2018; Month = 2; Day = 22; Timetable = "Purple"; Points = 6 },
    @{ Row = 14; Year = 2018; Month = 2; Day = 24; Timetable = "Purple"; Points = 6 },
    @{ Row = 15; Year = 2018; Month = 2; Day = 25; Timetable = "Purple"; Points = 6 }
)

foreach ($info in $rows) {
    $r = $info.Row

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dateCell.Value = Get-Date -Year $info.Year -Month $info.Month -Day $info.Day -Hour 0 -Minute 0 -Second 0

    $ws.Cells.Item($r, 2).Value = $info.Timetable
    $ws.Cells.Item($r, 4).Value = $info.Points
}

Write-Host "roster updated"
